# bun-rieu.xlsx - "Zutaten angepasst bun rieu" (ingredients adjusted)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tomatenmark: 100 -> 200 g
$ws.Range("C3").Value = 200

# Zwiebeln: 6 -> 5 Stück
$ws.Range("C6").Value = 5

# Eier (M) -> Eier (L), quantity 6 -> 4
$ws.Range("B7").Value = "Eier (L)"
$ws.Range("C7").Value = 4

# Wasser: 6 -> 1 l
$ws.Range("C8").Value = 1

# Oel: 4 -> 3, Teeloeffel -> Essloeffel
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "Esslöffel"

# Fruehlingszwiebeln: 2 -> 3 Stueck
$ws.Range("C14").Value = 3

# active cell moved from C8 to B8
$ws.Range("B8").Select()
